$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H19").Value = 5166.2
$ws.Range("I19").Value = 4939
$ws.Range("J19").Value = 5223
$ws.Range("K19").Value = 4939
$ws.Range("L19").Value = 5223
$ws.Range("M19").Value = -4764
$ws.Range("N19").Value = -5573
$ws.Range("H28").Value = 2016
$ws.Range("I28").Value = 1133
$ws.Range("K28").Value = 1133
$ws.Range("M28").Value = -648
$ws.Range("H33").Value = 5458488
$ws.Range("I33").Value = 5954692
$ws.Range("K33").Value = 5954692
$ws.Range("M33").Value = -5954463
$ws.Range("H40").Value = 6355
$ws.Range("I40").Value = 4993
$ws.Range("J40").Value = 7263
$ws.Range("K40").Value = 4993
$ws.Range("L40").Value = 7263
$ws.Range("M40").Value = -4818
$ws.Range("N40").Value = -7613
$ws.Range("H55").Value = 130.14706
$ws.Range("I55").Value = 138.8125
$ws.Range("J55").Value = 122.44444
$ws.Range("K55").Value = 138.8125
$ws.Range("L55").Value = 122.44444
$ws.Range("M55").Value = 75.1875
$ws.Range("N55").Value = -550.44444
$ws.Range("H57").Value = 54979
$ws.Range("J57").Value = 54979
$ws.Range("L57").Value = 164937
$ws.Range("N57").Value = -165935
$ws.Range("H76").Value = 4253.8823
$ws.Range("I76").Value = 4138.625
$ws.Range("K76").Value = 4138.625
$ws.Range("M76").Value = -3823.625
$ws.Range("H79").Value = 4253.8823
$ws.Range("I79").Value = 4138.625
$ws.Range("K79").Value = 4138.625
$ws.Range("M79").Value = -3046.625
$ws.Range("H98").Value = 2509.7693
$ws.Range("I98").Value = 2085.1
$ws.Range("J98").Value = 3925.3333
$ws.Range("K98").Value = 2085.1
$ws.Range("L98").Value = 3925.3333
$ws.Range("M98").Value = -587.0999999999999
$ws.Range("N98").Value = -6921.3333
$ws.Range("H107").Value = 696.8570999999999
$ws.Range("I107").Value = 768.9091
$ws.Range("J107").Value = 432.66666
$ws.Range("K107").Value = 768.9091
$ws.Range("L107").Value = 432.66666
$ws.Range("M107").Value = 1151.0909
$ws.Range("N107").Value = -4272.66666
$ws.Range("H113").Value = 7884.5713
$ws.Range("I113").Value = 7918.6
$ws.Range("J113").Value = 7799.5
$ws.Range("K113").Value = 7918.6
$ws.Range("L113").Value = 7799.5
$ws.Range("M113").Value = -4664.6
$ws.Range("N113").Value = -14307.5
$ws.Range("H122").Value = 2509.7693
$ws.Range("I122").Value = 2085.1
$ws.Range("J122").Value = 3925.3333
$ws.Range("K122").Value = 6255.299999999999
$ws.Range("L122").Value = 11775.9999
$ws.Range("M122").Value = -3805.299999999999
$ws.Range("N122").Value = -16675.9999
$ws.Range("H132").Value = 79113.84
$ws.Range("I132").Value = 85290
$ws.Range("K132").Value = 255870
$ws.Range("M132").Value = -253340
$ws.Range("H138").Value = 2925.122
$ws.Range("J138").Value = 3808.4644
$ws.Range("L138").Value = 11425.3932
$ws.Range("N138").Value = -21705.3932

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 8356
$ws.Range("I25").Value = 4468.6665
$ws.Range("K25").Value = 4468.6665
$ws.Range("M25").Value = -4066.6665
$ws.Range("H45").Value = 1267.25
$ws.Range("I45").Value = 1037.2
$ws.Range("K45").Value = 1037.2
$ws.Range("M45").Value = -660.2
$ws.Range("H110").Value = 4420.9614
$ws.Range("I110").Value = 4087.2727
$ws.Range("K110").Value = 4087.2727
$ws.Range("M110").Value = -2042.2727
$ws.Range("H122").Value = 2586.476
$ws.Range("I122").Value = 2471.7856
$ws.Range("K122").Value = 7415.3568
$ws.Range("M122").Value = -4965.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1275.6364
$ws.Range("I20").Value = 1127.3529
$ws.Range("J20").Value = 1779.8
$ws.Range("K20").Value = 1127.3529
$ws.Range("L20").Value = 1779.8
$ws.Range("M20").Value = -880.3529000000001
$ws.Range("N20").Value = -2273.8
$ws.Range("H25").Value = 6644.909
$ws.Range("I25").Value = 2008.6666
$ws.Range("J25").Value = 27508
$ws.Range("K25").Value = 2008.6666
$ws.Range("L25").Value = 27508
$ws.Range("M25").Value = -1773.6666
$ws.Range("N25").Value = -27978
$ws.Range("H100").Value = 13999.5
$ws.Range("J100").Value = 13999.5
$ws.Range("L100").Value = 13999.5
$ws.Range("N100").Value = -16163.5
$ws.Range("H130").Value = 99819.2
$ws.Range("J130").Value = 99819.2
$ws.Range("L130").Value = 99819.2
$ws.Range("N130").Value = -109859.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14927.789
$ws.Range("J31").Value = 3496.6667
$ws.Range("L31").Value = 3496.6667
$ws.Range("N31").Value = -4086.6667
$ws.Range("H34").Value = 14927.789
$ws.Range("J34").Value = 3496.6667
$ws.Range("L34").Value = 3496.6667
$ws.Range("N34").Value = -3900.6667
$ws.Range("H107").Value = 2231.4
$ws.Range("I107").Value = 854.9
$ws.Range("K107").Value = 854.9
$ws.Range("M107").Value = 1065.1
$ws.Range("H132").Value = 2438.5833
$ws.Range("I132").Value = 2486.75
$ws.Range("K132").Value = 7460.25
$ws.Range("M132").Value = -4930.25
$ws.Range("H134").Value = 5557288
$ws.Range("I134").Value = 1647.3334
$ws.Range("J134").Value = 16668569
$ws.Range("K134").Value = 4942.0002
$ws.Range("L134").Value = 50005707
$ws.Range("M134").Value = -2407.0002
$ws.Range("N134").Value = -50010777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1569925
$ws.Range("I4").Value = 991.4
$ws.Range("J4").Value = 5056444
$ws.Range("K4").Value = 2974.2
$ws.Range("L4").Value = 15169332
$ws.Range("M4").Value = -2862.2
$ws.Range("N4").Value = -15169556
$ws.Range("H17").Value = 172.125
$ws.Range("I17").Value = 183
$ws.Range("J17").Value = 165.6
$ws.Range("K17").Value = 549
$ws.Range("L17").Value = 496.8
$ws.Range("M17").Value = -380
$ws.Range("N17").Value = -834.8
$ws.Range("H64").Value = 14431.429
$ws.Range("I64").Value = 5167.3335
$ws.Range("K64").Value = 15502.0005
$ws.Range("M64").Value = -15232.0005
$ws.Range("H67").Value = 14431.429
$ws.Range("I67").Value = 5167.3335
$ws.Range("K67").Value = 15502.0005
$ws.Range("M67").Value = -14566.0005
$ws.Range("H116").Value = 104804.766
$ws.Range("I116").Value = 112788.5
$ws.Range("K116").Value = 338365.5
$ws.Range("M116").Value = -334923.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2521346.8
$ws.Range("I24").Value = 10000000
$ws.Range("K24").Value = 10000000
$ws.Range("M24").Value = -9999827

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6736.231
$ws.Range("I7").Value = 7453.222
$ws.Range("J7").Value = 5123
$ws.Range("K7").Value = 7453.222
$ws.Range("L7").Value = 5123
$ws.Range("M7").Value = -7341.222
$ws.Range("N7").Value = -5347
$ws.Range("H40").Value = 6684.8
$ws.Range("I40").Value = 3335
$ws.Range("J40").Value = 7522.25
$ws.Range("K40").Value = 3335
$ws.Range("L40").Value = 7522.25
$ws.Range("M40").Value = -3199
$ws.Range("N40").Value = -7794.25
$ws.Range("H46").Value = 1783.5714
$ws.Range("I46").Value = 1099.2
$ws.Range("K46").Value = 1099.2
$ws.Range("M46").Value = -911.2
$ws.Range("H55").Value = 2874.75
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 2833
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 2833
$ws.Range("M55").Value = -2827
$ws.Range("N55").Value = -3179
$ws.Range("H61").Value = 2285
$ws.Range("J61").Value = 4622.25
$ws.Range("L61").Value = 4622.25
$ws.Range("N61").Value = -5026.25
$ws.Range("H82").Value = 5641.385
$ws.Range("I82").Value = 4121.1816
$ws.Range("K82").Value = 4121.1816
$ws.Range("M82").Value = -3760.1816
$ws.Range("H85").Value = 5641.385
$ws.Range("I85").Value = 4121.1816
$ws.Range("K85").Value = 4121.1816
$ws.Range("M85").Value = -2873.1816
$ws.Range("H95").Value = 25322
$ws.Range("J95").Value = 25322
$ws.Range("L95").Value = 25322
$ws.Range("N95").Value = -30814
$ws.Range("H113").Value = 2285
$ws.Range("J113").Value = 4622.25
$ws.Range("L113").Value = 4622.25
$ws.Range("N113").Value = -8962.25
$ws.Range("H122").Value = 3594.25
$ws.Range("J122").Value = 3594.25
$ws.Range("L122").Value = 10782.75
$ws.Range("N122").Value = -15682.75
$ws.Range("H126").Value = 6736.231
$ws.Range("I126").Value = 7453.222
$ws.Range("J126").Value = 5123
$ws.Range("K126").Value = 22359.666
$ws.Range("L126").Value = 15369
$ws.Range("M126").Value = -19889.666
$ws.Range("N126").Value = -20309
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 58339.668
$ws.Range("I24").Value = 55000
$ws.Range("K24").Value = 55000
$ws.Range("M24").Value = -54770
$ws.Range("H122").Value = 2708.577
$ws.Range("I122").Value = 2866.4119
$ws.Range("J122").Value = 2410.4443
$ws.Range("K122").Value = 8599.235700000001
$ws.Range("L122").Value = 7231.3329
$ws.Range("M122").Value = -6149.235700000001
$ws.Range("N122").Value = -12131.3329

Write-Output "Applied all changes"
